$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Title: "Title: Life history trade-offs " -> the new climate-adaptation
#    title (ending with two trailing spaces, matching the target runs).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Title: Life history trade-offs ", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Title: Climate associated selection produces non-linear patterns of local adaptation in physiological trade-offs.  ",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Helper: insert a brand-new paragraph right after paragraph #afterIdx,
#    populating it with one-or-more runs of text (in order). Returns the
#    index of the freshly created paragraph.
# ---------------------------------------------------------------------------
function Insert-ParaAfter($afterIdx, [string[]]$runs) {
    $anchor = $d.Paragraphs.Item($afterIdx).Range
    $anchor.InsertParagraphAfter()
    $newIdx = $afterIdx + 1
    $pr = $d.Paragraphs.Item($newIdx).Range
    $pr.Collapse(1) # wdCollapseStart
    foreach ($run in $runs) {
        if ($run -ne "") {
            $pr.InsertAfter($run)
            $pr.Collapse(0) # wdCollapseEnd -- advance past the text just added
        }
    }
    return $newIdx
}

# ---------------------------------------------------------------------------
# 3. Build the new "Latitudinal variation ... Signatures of ... trade-offs"
#    block. It is inserted right after the blank paragraph that follows
#    "Introduction" (paragraph 10 in the original document), i.e. before the
#    "_GoBack" bookmark paragraph.
# ---------------------------------------------------------------------------
$introBlankIdx = 10
if ($d.Paragraphs.Item($introBlankIdx).Range.Text.Trim() -ne "") {
    throw "Unexpected document layout: paragraph $introBlankIdx is not blank"
}

$idx = $introBlankIdx
$idx = Insert-ParaAfter $idx @("Latitudinal variation is lacking in china, so that makes comparisons for parallel evolution sparse. ")
$idx = Insert-ParaAfter $idx @()
$idx = Insert-ParaAfter $idx @()
$idx = Insert-ParaAfter $idx @()
$idx = Insert-ParaAfter $idx @("Materials and methods")
$idx = Insert-ParaAfter $idx @()
$idx = Insert-ParaAfter $idx @()
$idx = Insert-ParaAfter $idx @()
$idx = Insert-ParaAfter $idx @("Results")
$idx = Insert-ParaAfter $idx @()
$idx = Insert-ParaAfter $idx @("Trade-offs in ", "physiology")
$idx = Insert-ParaAfter $idx @()
$idx = Insert-ParaAfter $idx @("Patterns of climate across latitude")
$idx = Insert-ParaAfter $idx @()
$listIdx = Insert-ParaAfter $idx @("How different climate pcs correlate with latitude")
$idx = $listIdx
$idx = Insert-ParaAfter $idx @()
$idx = Insert-ParaAfter $idx @("Signatures of", " local adaptation in ", "physiological ", "trade-offs")

$afterBlockIdx = $idx

# ---------------------------------------------------------------------------
# 4. Turn the "How different climate pcs correlate with latitude" paragraph
#    into a bulleted list item (List Paragraph style + numId 1), matching
#    the numbering.xml / styles.xml additions in the target document.
# ---------------------------------------------------------------------------
$listRange = $d.Paragraphs.Item($listIdx).Range
$listRange.Style = "List Paragraph"
$listTemplate = $word.ListGalleries.Item(1).ListTemplates.Item(1)
$listRange.ListFormat.ApplyListTemplateWithLevel($listTemplate)
try {
    $listStyle = $d.Styles.Item("List Paragraph")
    $listStyle.Priority = 34
    $listStyle.ParagraphFormat.LeftIndent = 36
    $listStyle.NoSpaceBetweenParagraphsOfSameStyle = $true
} catch {
}

# ---------------------------------------------------------------------------
# 5. The original "Materials and methods" / blank / "Results" paragraphs
#    (which used to sit right after the "_GoBack" bookmark paragraph) have
#    now been superseded by the copies inserted above, so remove the old
#    trio. They immediately follow the bookmark paragraph, which itself
#    immediately follows the newly inserted block.
# ---------------------------------------------------------------------------
$bookmarkIdx = $afterBlockIdx + 1
$oldMaterialsIdx = $bookmarkIdx + 1
$oldBlankIdx = $bookmarkIdx + 2
$oldResultsIdx = $bookmarkIdx + 3

if ($d.Paragraphs.Item($oldMaterialsIdx).Range.Text.Trim() -ne "Materials and methods" -or `
    $d.Paragraphs.Item($oldResultsIdx).Range.Text.Trim() -ne "Results") {
    throw "Unexpected document layout: could not locate the old Materials/Results paragraphs"
}

$delStart = $d.Paragraphs.Item($oldMaterialsIdx).Range.Start
$delEnd = $d.Paragraphs.Item($oldResultsIdx).Range.End
$d.Range($delStart, $delEnd).Delete()

Write-Host "Done. Final paragraph count:" $d.Paragraphs.Count
